# Commit: "Automatic update of files." -- the canonical XML diff shows that
# rows 5 and 7 swap their entire contents, and rows 6 and 8 swap their
# entire contents (every column A:AY), on the single worksheet "Artfynd".
#
# We implement this as a generic row-swap: snapshot every cell in the
# affected rows first (so later writes don't clobber values we still need
# to read), then write the swapped rows back, being careful to:
#   - keep numeric columns numeric
#   - keep boolean columns boolean
#   - keep everything else as literal text (so values that merely look
#     like numbers/dates, e.g. "1" or "2020-06-12", are not silently
#     reinterpreted by Excel's automatic type detection)
#   - fully clear destination cells whose source was empty, so columns
#     that should be completely absent for a row stay absent

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$maxCol = 51  # column AY

# Columns that hold real numbers: A, B, E, Q, R, S
$numericCols = @(1, 2, 5, 17, 18, 19)
# Columns that hold booleans: AD, AE, AG
$boolCols = @(30, 31, 33)

# Rows participating in the swap, and where each destination row's data
# comes from.
$swapMap = @{5 = 7; 6 = 8; 7 = 5; 8 = 6}
$rowsInvolved = @(5, 6, 7, 8)

function Test-LooksSpecial($val) {
    # Values that Excel's auto-detection would reinterpret as a number,
    # date, or boolean if written without an explicit text format.
    if ($val -match '^\s*[-+]?\d+(\.\d+)?\s*$') { return $true }
    if ($val -match '^\s*\d{1,4}[-/]\d{1,2}[-/]\d{1,4}\s*$') { return $true }
    if ($val -match '^(?i:true|false)$') { return $true }
    return $false
}

# Step 1: snapshot current values of every involved row before any writes.
# For numeric/boolean columns we read .Value2 (typed, full precision). For
# every other column we read .Text instead of .Value2: this engine's COM
# layer silently coerces digit-only *text* cells (e.g. the literal string
# "1" in column I) to a Double when read through .Value2, which would lose
# the fact that the source cell is really text. .Text always returns the
# literal display string, sidestepping that coercion.
$data = @{}
foreach ($r in $rowsInvolved) {
    $rowVals = @()
    for ($col = 1; $col -le $maxCol; $col++) {
        $cell = $ws.Cells.Item($r, $col)
        if ($numericCols -contains $col -or $boolCols -contains $col) {
            $rowVals += $cell.Value2
        } else {
            $rowVals += $cell.Text
        }
    }
    $data[$r] = $rowVals
}

# Step 2: write each destination row from its snapshot source row.
foreach ($destRow in $rowsInvolved) {
    $srcRow = $swapMap[$destRow]
    $vals = $data[$srcRow]
    for ($col = 1; $col -le $maxCol; $col++) {
        $cell = $ws.Cells.Item($destRow, $col)
        $val = $vals[$col - 1]

        if ($numericCols -contains $col -or $boolCols -contains $col) {
            $cell.Value2 = $val
        } elseif ($null -eq $val -or $val -eq "") {
            # Source was empty/absent -> destination becomes empty/absent.
            $cell.Value2 = $null
        } else {
            $text = [string]$val
            if (Test-LooksSpecial $text) {
                $cell.NumberFormat = "@"
            }
            $cell.Value2 = $text
        }
    }
}
